# The authored change deletes the original "index" column (the first
# data column, containing 0,1,2,...,99) from the worksheet. Excel's
# native "delete entire column" shifts every remaining column one slot
# to the left, re-maps cell styles/borders accordingly, and shrinks the
# used range from A1:G101 down to A1:F101 - exactly matching the target
# diff (headers move from B1:F1 into A1:E1, the stray summary value
# moves from G1 into F1, and each data row's values shift left by one
# column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()
